# Fix typo in the D100 magic item table on "Feuil1": the last range should
# read "96-100" (not "96-00") since the table is a d100 roll table.
# Also restores the last-used cell selection on that sheet to F21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

$ws.Range("A23").Value = "96-100 "

$ws.Range("F21").Select()
